$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Femacal de La Calera - Ciboulette".
# Insert a fresh row at 125 (pushing the existing rows 125-353 down to
# 126-354, so the former last row, 353, becomes 354) and populate it with
# the new record's data.
$ws.Rows.Item(125).Insert()

$ws.Range("A125").Value = 3
$ws.Range("B125").Value = "Femacal de La Calera"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44791
$ws.Range("E125").Value = 5
$ws.Range("F125").Value = 100112039
$ws.Range("G125").Value = "Ciboulette"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 160
$ws.Range("K125").Value = 1500
$ws.Range("L125").Value = 1500
$ws.Range("M125").Value = 1500
$ws.Range("N125").Value = '$/docena de atados'
$ws.Range("O125").Value = "Provincia de Quillota"
$ws.Range("P125").Value = 500
$ws.Range("Q125").Value = 3
$ws.Range("R125").Value = "Hortaliza"
